{"js": "// Remove the trailing \"Ver no Jupiter...\" and \"\u00a9 2020 ...\" boilerplate\n// paragraphs (and the blank paragraph that precedes them), which the\n// site generator no longer emits after a rebuild.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst targetSnippets = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\nconst items = paragraphs.items;\n\n// Find the paragraphs that hold the two known text snippets.\nlet firstIndex = -1;\nlet lastIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  const text = items[i].text;\n  if (targetSnippets.some((snippet) => text.indexOf(snippet) !== -1)) {\n    if (firstIndex === -1) firstIndex = i;\n    lastIndex = i;\n  }\n}\n\nif (firstIndex !== -1) {\n  // Also drop the single empty paragraph immediately before the first\n  // matched paragraph (it only existed as a spacer for the boilerplate\n  // block that is being removed).\n  let startIndex = firstIndex;\n  if (startIndex > 0 && items[startIndex - 1].text === \"\") {\n    startIndex -= 1;\n  }\n\n  for (let i = lastIndex; i >= startIndex; i--) {\n    items[i].delete();\n  }\n\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" and \"(c) 2020 ...\" boilerplate\n# paragraphs (and the blank spacer paragraph right before them), which the\n# site generator no longer emits after a rebuild.\n\n$d = $word.ActiveDocument\n\n$snippet1 = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$snippet2 = \"Contact: luizeleno@usp.br\"\n\n$firstIndex = -1\n$lastIndex = -1\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n  $t = $d.Paragraphs.Item($i).Range.Text\n  if ($t.Contains($snippet1) -or $t.Contains($snippet2)) {\n    if ($firstIndex -eq -1) { $firstIndex = $i }\n    $lastIndex = $i\n  }\n}\n\nif ($firstIndex -ne -1) {\n  $startIndex = $firstIndex\n  # Also remove the single empty spacer paragraph immediately before the\n  # matched block (its own text is just the paragraph mark, length 1).\n  if ($startIndex -gt 1) {\n    $prevText = $d.Paragraphs.Item($startIndex - 1).Range.Text\n    if ($prevText.Trim().Length -eq 0) {\n      $startIndex = $startIndex - 1\n    }\n  }\n\n  for ($i = $lastIndex; $i -ge $startIndex; $i--) {\n    $d.Paragraphs.Item($i).Range.Delete()\n  }\n}\n"}
